# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet: bumps the data timestamp, re-ranks a
# handful of countries whose case counts crossed each other (so their name
# cells swap), and writes the latest totals/actives/recovered/etc. for the
# rows whose statistics changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 13:34"

# Country name swaps/shifts caused by re-ranking (column A text)
$ws.Range("A114").Value = "San Marino"
$ws.Range("A115").Value = "Georgia"
$ws.Range("A138").Value = "Zambia"
$ws.Range("A139").Value = "Benin"
$ws.Range("A140").Value = "Cabo Verde"
$ws.Range("A141").Value = "Etiopia"
$ws.Range("A142").Value = "Santo Tome y Principe"
$ws.Range("A143").Value = "Liberia"
$ws.Range("A144").Value = "Madagascar"
$ws.Range("A145").Value = "Islas Feroe"
$ws.Range("A146").Value = "Martinica"
$ws.Range("A147").Value = "Birmania"

# Updated case statistics (columns B-H)
$ws.Range("B52").Value = 7623
$ws.Range("C52").Value = 415
$ws.Range("D52").Value = 2622
$ws.Range("E52").Value = 4952
$ws.Range("F52").Value = 95
$ws.Range("G52").Value = 2
$ws.Range("H52").Value = 49

$ws.Range("B77").Value = 2090
$ws.Range("C77").Value = 20
$ws.Range("D77").Value = 1059
$ws.Range("E77").Value = 929
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 102

$ws.Range("B103").Value = 844
$ws.Range("C103").Value = 9
$ws.Range("E103").Value = 580

$ws.Range("B104").Value = 809
$ws.Range("C104").Value = 13
$ws.Range("D104").Value = 234
$ws.Range("E104").Value = 549

$ws.Range("B114").Value = 637
$ws.Range("C114").Value = 14
$ws.Range("D114").Value = 126
$ws.Range("E114").Value = 470
$ws.Range("F114").Value = 3
$ws.Range("H114").Value = 41

$ws.Range("B115").Value = 626
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 297
$ws.Range("E115").Value = 319
$ws.Range("F115").Value = 6
$ws.Range("H115").Value = 10

$ws.Range("B124").Value = 490
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 427
$ws.Range("E124").Value = 58

$ws.Range("B138").Value = 252
$ws.Range("C138").Value = 85
$ws.Range("D138").Value = 112
$ws.Range("E138").Value = 133
$ws.Range("F138").Value = 1
$ws.Range("G138").Value = 3
$ws.Range("H138").Value = 7

$ws.Range("B139").Value = 242
$ws.Range("D139").Value = 62
$ws.Range("E139").Value = 178

$ws.Range("B140").Value = 230
$ws.Range("C140").Value = 0
$ws.Range("D140").Value = 44
$ws.Range("E140").Value = 184
$ws.Range("F140").Value = 0
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 2

$ws.Range("B141").Value = 210
$ws.Range("C141").Value = 16
$ws.Range("D141").Value = 97
$ws.Range("E141").Value = 108
$ws.Range("F141").Value = 1
$ws.Range("G141").Value = 1

$ws.Range("B142").Value = 208
$ws.Range("D142").Value = 4
$ws.Range("E142").Value = 199
$ws.Range("H142").Value = 5

$ws.Range("B143").Value = 199
$ws.Range("D143").Value = 79
$ws.Range("E143").Value = 100
$ws.Range("F143").Value = 0
$ws.Range("H143").Value = 20

$ws.Range("B144").Value = 193
$ws.Range("D144").Value = 101
$ws.Range("E144").Value = 92
$ws.Range("F144").Value = 1

$ws.Range("B145").Value = 187
$ws.Range("D145").Value = 187
$ws.Range("E145").Value = 0
$ws.Range("F145").Value = 0
$ws.Range("H145").Value = 0

$ws.Range("B146").Value = 186
$ws.Range("D146").Value = 83
$ws.Range("E146").Value = 89
$ws.Range("F146").Value = 3
$ws.Range("H146").Value = 14

$ws.Range("B147").Value = 177
$ws.Range("D147").Value = 67
$ws.Range("E147").Value = 104
$ws.Range("F147").Value = 0
$ws.Range("H147").Value = 6

$ws.Range("D161").Value = 31
$ws.Range("E161").Value = 78

$ws.Range("D184").Value = 19
$ws.Range("E184").Value = 3
